$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns Q, R, S (old "per person US/WA" and "per vehicle WA" mileage columns)
# This shifts old column T (vehicle age) into Q, and keeps P (will be updated below)
$ws.Range("Q1:S1").EntireColumn.Delete()

# Update header text for column P (now represents avg annual mileage per vehicle U.S.)
$ws.Range("P1").Value = "average annual mileage per vehicle (U.S.) [mi]"

# Set flat custom mileage assumption of 13000 for all rows
$ws.Range("P2:P7").Value = 13000

# Update view state: selection
$ws.Range("P11").Select()

Write-Output "edit applied"
